$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header summary figures ---
# VALOR MORA total
$ws.Range("E11").Value2 = 225411

# Cant. Trabajadores / Cant. Periodos
$ws.Range("C13").Value2 = 3
$ws.Range("F13").Value2 = 7

# --- Update Salario Basico (column G) for the existing detail rows 16-21 ---
$ws.Range("G16").Value2 = 877803
$ws.Range("G17").Value2 = 877803
$ws.Range("G18").Value2 = 877803
$ws.Range("G19").Value2 = 877803
$ws.Range("G20").Value2 = 877803
$ws.Range("G21").Value2 = 1000000

# --- Insert a new detail row (row 22) for a new worker, shifting the
#     trailing signature rows (old 26/27 -> 27/28) down by one ---
$ws.Range("B22:J22").Insert(-4121) # xlShiftDown

# Clone the formatting of row 21 (borders/fill/font) onto the new row 22
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122) # xlPasteAll
$excel.CutCopyMode = 0

# Fill in the new worker's data
$ws.Range("B22").Value2 = "CC"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value2 = "1083007268"

$ws.Range("D22").Value2 = "SUSARAI PATRICIA PERTUZ BARROS"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "2508"

$ws.Range("F22").Value2 = 24674
$ws.Range("G22").Value2 = 1423500
